$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text formatting so values
# like "309.51" or "17.50" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '39.949.92'
$ws.Range("E2").Value = '  -2.90%  '
$ws.Range("D3").Value = '2.340.81'
$ws.Range("E3").Value = '  -3.66%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '309.51'
$ws.Range("E5").Value = '  -2.81%  '
$ws.Range("D6").Value = '84.23'
$ws.Range("E6").Value = '  -5.85%  '
$ws.Range("E7").Value = '  -2.85%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '0.479'
$ws.Range("D10").Value = '0.0803'
$ws.Range("E10").Value = '  -4.02%  '
$ws.Range("D11").Value = '29.89'
$ws.Range("E11").Value = '  -6.78%  '
$ws.Range("E12").Value = '  +0.66%  '
$ws.Range("D13").Value = '2.706.81'
$ws.Range("E13").Value = '  -3.49%  '
$ws.Range("D14").Value = '6.39'
$ws.Range("E14").Value = '  -5.17%  '
$ws.Range("D15").Value = '14.75'
$ws.Range("E15").Value = '  -6.44%  '
$ws.Range("D16").Value = '2.357.65'
$ws.Range("E16").Value = '  -2.99%  '
$ws.Range("D17").Value = '0.755'
$ws.Range("E17").Value = '  -2.71%  '
$ws.Range("D18").Value = '39.979.58'
$ws.Range("E18").Value = '  -2.71%  '
$ws.Range("D19").Value = '0.0₃0898'
$ws.Range("E19").Value = '  -3.14%  '
$ws.Range("D20").Value = '6.05'
$ws.Range("E20").Value = '  -3.83%  '
$ws.Range("D21").Value = '67.90'
$ws.Range("E21").Value = '  -5.80%  '
$ws.Range("D22").Value = '10.53'
$ws.Range("E22").Value = '  -5.00%  '
$ws.Range("D23").Value = '234.79'
$ws.Range("E23").Value = '  -0.16%  '
$ws.Range("E24").Value = '  -5.93%  '
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("D27").Value = '23.48'
$ws.Range("E27").Value = '  -2.91%  '
$ws.Range("D28").Value = '2.19'
$ws.Range("E28").Value = '  -1.64%  '
$ws.Range("D29").Value = '9.21'
$ws.Range("E29").Value = '  -4.25%  '
$ws.Range("D30").Value = '34.61'
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("D31").Value = '153.20'
$ws.Range("E31").Value = '  -2.85%  '
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("D33").Value = '5.06'
$ws.Range("E33").Value = '  -3.98%  '
$ws.Range("D34").Value = '2.48'
$ws.Range("E34").Value = '  -1.13%  '
$ws.Range("D35").Value = '0.0714'
$ws.Range("E35").Value = '  -4.27%  '
$ws.Range("E36").Value = '  -1.01%  '
$ws.Range("E37").Value = '  -6.24%  '
$ws.Range("D38").Value = '0.0985'
$ws.Range("E38").Value = '  -2.12%  '
$ws.Range("D39").Value = '15.49'
$ws.Range("E39").Value = '  -8.23%  '
$ws.Range("E40").Value = '  -4.07%  '
$ws.Range("D41").Value = '3.83'
$ws.Range("E41").Value = '  -1.73%  '
$ws.Range("D42").Value = '1.972.07'
$ws.Range("E42").Value = '  -1.20%  '
$ws.Range("D43").Value = '2.26'
$ws.Range("E43").Value = '  -0.50%  '
$ws.Range("D44").Value = '0.0264'
$ws.Range("E44").Value = '  -4.46%  '
$ws.Range("D45").Value = '17.50'
$ws.Range("E45").Value = '  -5.87%  '
$ws.Range("E46").Value = '  -1.70%  '
$ws.Range("D47").Value = '2.67'
$ws.Range("E47").Value = '  -7.91%  '
$ws.Range("D48").Value = '2.565.46'
$ws.Range("E48").Value = '  -3.80%  '
$ws.Range("D49").Value = '92.49'
$ws.Range("E49").Value = '  -2.40%  '
$ws.Range("D50").Value = '69.93'
$ws.Range("E50").Value = '  -4.92%  '
$ws.Range("D51").Value = '49.70'
$ws.Range("E51").Value = '  -4.08%  '
